# The presentation originally uses the "Integral" theme (ppt/theme/theme1.xml,
# linked from the slide master) together with a second, unused "Office Theme"
# theme part (ppt/theme/theme2.xml, linked only from the notes master). The
# authored change swaps the two themes' contents: the presentation's live
# theme becomes the default "Office Theme" colors (and the spare theme part
# ends up holding the former "Integral" palette).
#
# The actually-rendered, in-use theme is the one hanging off the slide
# master (Design.SlideMaster.Theme) -- that's ppt/theme/theme1.xml, so we
# repoint every themed color there to the stock "Office" palette.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

function PackRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# msoThemeColorDark1..msoThemeColorFollowedHyperlink (indices 1-12), using
# the stock Office theme palette.
$themeColors.Item(1).RGB  = PackRGB 0x00 0x00 0x00   # Dark 1
$themeColors.Item(2).RGB  = PackRGB 0xFF 0xFF 0xFF   # Light 1
$themeColors.Item(3).RGB  = PackRGB 0x44 0x54 0x6A   # Dark 2
$themeColors.Item(4).RGB  = PackRGB 0xE7 0xE6 0xE6   # Light 2
$themeColors.Item(5).RGB  = PackRGB 0x5B 0x9B 0xD5   # Accent 1
$themeColors.Item(6).RGB  = PackRGB 0xED 0x7D 0x31   # Accent 2
$themeColors.Item(7).RGB  = PackRGB 0xA5 0xA5 0xA5   # Accent 3
$themeColors.Item(8).RGB  = PackRGB 0xFF 0xC0 0x00   # Accent 4
$themeColors.Item(9).RGB  = PackRGB 0x44 0x72 0xC4   # Accent 5
$themeColors.Item(10).RGB = PackRGB 0x70 0xAD 0x47   # Accent 6
$themeColors.Item(11).RGB = PackRGB 0x05 0x63 0xC1   # Hyperlink
$themeColors.Item(12).RGB = PackRGB 0x95 0x4F 0x72   # Followed Hyperlink
